$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column so numeric-looking values stay as text
# (matches original inlineStr cell type), then restore default formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "71.093.78"
$ws.Range("E2").Value = "  +3.07%  "

$ws.Range("D3").Value = "3.815.16"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "705.56"
$ws.Range("E5").Value = "  +11.90%  "

$ws.Range("D6").Value = "173.39"
$ws.Range("E6").Value = "  +4.77%  "

$ws.Range("D7").Value = "3.814.21"
$ws.Range("E7").Value = "  +1.24%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("E10").Value = "  +3.69%  "

$ws.Range("E11").Value = "  +10.02%  "

$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  +9.32%  "

$ws.Range("D14").Value = "36.37"
$ws.Range("E14").Value = "  +4.48%  "

$ws.Range("D15").Value = "4.458.58"
$ws.Range("E15").Value = "  +1.21%  "

$ws.Range("D16").Value = "3.808.01"
$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").Value = "71.119.92"
$ws.Range("E17").Value = "  +3.12%  "

$ws.Range("D18").Value = "18.03"
$ws.Range("E18").Value = "  +2.06%  "

$ws.Range("D19").Value = "7.26"
$ws.Range("E19").Value = "  +3.84%  "

$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("D21").Value = "11.27"
$ws.Range("E21").Value = "  +18.64%  "

$ws.Range("D22").Value = "483.80"
$ws.Range("E22").Value = "  +3.60%  "

$ws.Range("D23").Value = "0.718"
$ws.Range("E23").Value = "  +2.27%  "

$ws.Range("D24").Value = "83.97"
$ws.Range("E24").Value = "  +2.23%  "

$ws.Range("E25").Value = "  +2.45%  "

$ws.Range("D26").Value = "12.49"
$ws.Range("E26").Value = "  +3.45%  "

$ws.Range("D27").Value = "10.63"
$ws.Range("E27").Value = "  +5.10%  "

$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  +3.42%  "

$ws.Range("D29").Value = "3.967.22"
$ws.Range("E29").Value = "  +1.16%  "

$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").Value = "3.09"
$ws.Range("E31").Value = "  +15.60%  "

$ws.Range("D32").Value = "2.31"
$ws.Range("E32").Value = "  +2.16%  "

$ws.Range("D33").Value = "7.56"
$ws.Range("E33").Value = "  +6.64%  "

$ws.Range("D34").Value = "29.61"
$ws.Range("E34").Value = "  +4.26%  "

$ws.Range("E35").Value = "  +2.44%  "

$ws.Range("D36").Value = "9.26"
$ws.Range("E36").Value = "  +4.16%  "

$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").Value = "3.766.34"
$ws.Range("E38").Value = "  +1.15%  "

$ws.Range("E39").Value = "  +3.80%  "

$ws.Range("D40").Value = "3.51"
$ws.Range("E40").Value = "  +7.48%  "

$ws.Range("D41").Value = "6.01"
$ws.Range("E41").Value = "  +4.14%  "

$ws.Range("D42").Value = "0.000341"
$ws.Range("E42").Value = "  +30.20%  "

$ws.Range("D43").Value = "2.25"
$ws.Range("E43").Value = "  +13.67%  "

$ws.Range("D44").Value = "0.976"
$ws.Range("E44").Value = "  +1.46%  "

$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D47").Value = "45.39"
$ws.Range("E47").Value = "  +4.12%  "

$ws.Range("D48").Value = "161.08"
$ws.Range("E48").Value = "  +3.07%  "

$ws.Range("D49").Value = "49.42"
$ws.Range("E49").Value = "  +5.34%  "

$ws.Range("D50").Value = "1.42"
$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("E51").Value = "  +2.85%  "

$ws.Range("D2:D51").ClearFormats()